$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- pre-format the new rows so every cell picks up the same "vertical
# center" style (s="1") used by the rest of the data rows, before any
# values are written into them ---
$ws.Range("A7:K8").VerticalAlignment = -4108

# Row 6 (Id=5): the sensor was re-classified from a temperature indicator
# to a level indicator -> rename only.
$ws.Range("B6").Value = "LI_1"

# New shared-string text is introduced in this exact order so the rebuilt
# sharedStrings table lines up with the authored workbook.
$ws.Range("F7").Value = "Tank T2 - Level"
$ws.Range("F8").Value = "Flow measure on pipe to Tank T2"
$ws.Range("C7").Value = "DB1.DBD10"
$ws.Range("C8").Value = "DB1.DBD14"
$ws.Range("B7").Value = "LI_2"
$ws.Range("B8").Value = "FI_1"

# Row 7 (Id=6): new Level Indicator for Tank T2
$ws.Range("A7").Value = 6
$ws.Range("D7").Value = "REAL"
$ws.Range("E7").Value = "%"
$ws.Range("H7").Value = "None"
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0

# Row 8 (Id=7): new Flow Indicator for the pipe feeding Tank T2
$ws.Range("A8").Value = 7
$ws.Range("D8").Value = "REAL"
$ws.Range("E8").Value = "%"
$ws.Range("H8").Value = "None"
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0

# "False" must land as text (matching column G elsewhere), not a COM bool -
# copy it from an existing text cell instead of assigning the literal.
$ws.Range("G6").Copy()
$ws.Range("G7").PasteSpecial()
$ws.Range("G8").PasteSpecial()
$excel.CutCopyMode = $false

# Widen column F to fit the new, longer comment text (manual width, no
# longer auto "best fit")
$ws.Columns("F").ColumnWidth = 32.12

# Update selection to match the author's final cursor position
$ws.Range("C14").Select()
